$wb = $excel.ActiveWorkbook

$ws2025 = $wb.Worksheets.Item("2025")
$ws2025.Range("A2").Value = 0
$ws2025.Range("E2").Value = 297405.7209542279
$ws2025.Range("G2").Value = 80959.25712661858
$ws2025.Range("I2").Value = 142007.2657838
$ws2025.Range("L2").Value = 533802.9012525
$ws2025.Range("M2").Value = 106583.87169815
$ws2025.Range("N2").Value = 75545.14729476102
$ws2025.Range("O2").Value = 70974.1014704362

$ws2030 = $wb.Worksheets.Item("2030")
$ws2030.Range("A2").Value = 0
$ws2030.Range("B2").Value = 32088.98512992542
$ws2030.Range("E2").Value = 174532.7286598571
$ws2030.Range("I2").Value = 96226.01827323649
$ws2030.Range("L2").Value = 224920.9130748752
$ws2030.Range("M2").Value = 58508.23124971001
$ws2030.Range("N2").Value = 17274.62728995873
$ws2030.Range("O2").Value = 12264.56043766227

$ws2035 = $wb.Worksheets.Item("2035")
$ws2035.Range("A2").Value = 15535.6866473945
$ws2035.Range("B2").Value = 21697.4750396202
$ws2035.Range("E2").Value = 133319.8868932023
$ws2035.Range("I2").Value = 177127.1711388644
$ws2035.Range("L2").Value = 0
$ws2035.Range("M2").Value = 55808.89021025997
$ws2035.Range("N2").Value = 43148.97806764329
$ws2035.Range("O2").Value = 44476.65803456417
